$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a soft hyphen as its own run immediately before "Marius Olariu"
#    (first paragraph of the document).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.Collapse(1)
$r1.InsertBefore([char]31)          # Chr(31) == Word's internal soft-hyphen marker
$softHyphenRange = $d.Range(0, 1)
# Toggle Bold on/off (no visible effect) so this 1-character range keeps its
# own run instead of being silently re-merged into the following run.
$softHyphenRange.Bold = 1
$softHyphenRange.Bold = 0

# ---------------------------------------------------------------------------
# 2. Collapse the "...JavaScript, " / "and " / "React. ..." runs into one.
# ---------------------------------------------------------------------------
$find = "I am proficient in a diverse set of technologies including Java, AWS, Spring, and SQL, with a primary focus on back-end development but also hands-on experience with frontend technologies like JavaScript, and React. I am enthusiastic about learning new things, both technical and business-related."
$d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $find, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Collapse the "As a Backend Software Engineer..." runs into one.
# ---------------------------------------------------------------------------
$find = "As a Backend Software Engineer, I developed microservices with Java and Spring Boot and implemented initiatives that improved the pull request process, standardization across microservices in the department. I took the initiative in the management of the AWS infrastructure and gave talks on different technical subjects. As an engineer, I took responsibility and led different tech initiatives which has contributed to my promotion as team lead. My team was formed out of 5 developers, 1 Software Tester and 1 Scrum Master.  Additionally, my recognition as Top Talent, a distinction shared by only four individuals, acknowledged my technical and leadership skills."
$d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $find, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Fix the IQVIA employment date: "(Feb 2023 - Oct 2023)" -> "(Feb 2022 - Oct 2022)"
#    Reproduce the 5-run split seen in the target (the two "3"->"2" edits plus
#    the forced run boundaries around them), not just a plain merge.
# ---------------------------------------------------------------------------
$dateRange = $d.Content
$dateRange.Find.Execute("(Feb 2023 " + [char]0x2013 + " Oct 2023)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $dateRange.Start

# Fix the two digits first (same length edit, so offsets below stay valid).
$d.Range($start + 8, $start + 9).Text = "2"
$d.Range($start + 19, $start + 20).Text = "2"

# Force the run boundaries: "(Feb 202" | "2 " | "<EN DASH> Oct 202" | "2" | ")"
$segments = @(@(0, 8), @(8, 10), @(10, 19), @(19, 20), @(20, 21))
foreach ($seg in $segments) {
    $segRange = $d.Range($start + $seg[0], $start + $seg[1])
    $segRange.Bold = 1
    $segRange.Bold = 0
}

# ---------------------------------------------------------------------------
# 5. Collapse "...Product Owner" / "s" / " I delivered..." runs into one.
# ---------------------------------------------------------------------------
$find = "Full stack development for a Phase I eSource and clinical trial automation system which is used by the biggest pharmaceutical companies in the world. In collaboration with Project Managers, Biomedical Engineers and Product Owners I delivered for our clients high quality software, and I provided speedy troubleshooting, bug fixing and technical support. "
$d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $find, 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Collapse "...relevant " / "to" / " my work)..." runs into one (Diffblue).
# ---------------------------------------------------------------------------
$find = "I developed clean and performant software in Java for Diffblue Cover, an AI-based tool that automatically generates (unit) tests. I embraced continuous learning (i.e. each Sprint I studied a new topic relevant to my work) and took up new responsibilities (e.g. lead stand-ups and Sprint Retrospectives). Moreover, I set up and maintained an online environment where the product could be demoed. "
$d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $find, 2) | Out-Null

# ---------------------------------------------------------------------------
# 7. Collapse "...for Jira " / "which" / " a light-weight CRM..." runs into one.
# ---------------------------------------------------------------------------
$find = "I worked on the architecture and implementation of a cloud add-on for Jira which a light-weight CRM.  Also, I supported a colleague to transition from a Technical Consultant position to a Developer one through mentoring. In parallel to my work engagement, I completed my MSc dissertation."
$d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $find, 2) | Out-Null

# ---------------------------------------------------------------------------
# 8. Collapse the SFC scholarship sentence's runs ("which", "master'" + "s") into one.
# ---------------------------------------------------------------------------
$find = ". A scholarship awarded by the Scottish Funding Council (SFC) which covered the tuition fee for my master" + [char]0x2019 + "s degree programme. "
$d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $find, 2) | Out-Null
